# Updated cryptos list on Sat Apr  8 16:57:50 UTC 2023 with GitHub Actions
#
# The sheet lists ~50 cryptocurrencies with columns:
#   B = Coin name, C = coinranking.com link, D = Price, E = 1h volume/change.
# This refresh shifts several coins up/down the ranking (so B/C/D/E move as a
# block to a different row) and updates Price/Volume text for every row.
#
# D/E are plain text cells (e.g. "28.179.47", "  +0.74%  "), not numbers - the
# site renders thousands-separated/placeholder price strings that are not
# valid numeric literals. Assigning them straight to Range.Value is safe for
# most of them (Excel can't parse "28.179.47" as a number because of the
# double dot, and the volume strings have spaces/%), but plain single-dot
# decimals like "1.010" or "91.10" would be auto-coerced to the number
# 1.01 / 91.1 and lose the trailing zero/text type. For those we prefix the
# value with a leading apostrophe, exactly like typing '1.010 into Excel,
# which keeps it as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.179.47'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '1.879.34'
$ws.Range("E3").Value = '  +1.18%  '

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = '  +0.55%  '

$ws.Range("D5").Value = "'314.53"
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").Value = "'0.5121"
$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("D8").Value = "'0.3905"
$ws.Range("E8").Value = '  +1.87%  '

$ws.Range("D9").Value = "'0.08363"
$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").Value = "'41.64"
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").Value = "'6.229"
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'20.56"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.867.22'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("D15").Value = "'7.256"
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").Value = "'1.011"
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").Value = "'91.10"
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").Value = "'0.06682"
$ws.Range("E19").Value = '  +0.44%  '

$ws.Range("D20").Value = "'17.74"
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").Value = "'6.016"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").Value = '28.204.00'
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("D24").Value = "'11.10"
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").Value = "'2.246"
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.102.83'
$ws.Range("E26").Value = '  +1.51%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'159.82"
$ws.Range("E27").Value = '  +1.53%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.466"
$ws.Range("E28").Value = '  -1.88%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = "'20.76"
$ws.Range("E29").Value = '  +1.43%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = "'126.27"
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.1055"
$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'1.041"
$ws.Range("E32").Value = '  +0.95%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'5.846"
$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = "'3.607"
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = "'9.524"
$ws.Range("E35").Value = '  +1.58%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = "'0.02440"
$ws.Range("E36").Value = '  +1.23%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.06558"
$ws.Range("E37").Value = '  +0.96%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2211"
$ws.Range("E38").Value = '  +1.85%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = "'1.195"
$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6466"
$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.246"
$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").Value = "'4.984"
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'11.19"
$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'0.6076"
$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'12.98"
$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'3.688"
$ws.Range("E46").Value = '  +1.04%  '

$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'1.278"
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'2.011"
$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = "'1.233"
$ws.Range("E49").Value = '  +2.36%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'120.90"
$ws.Range("E50").Value = '  +0.81%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.06906"
$ws.Range("E51").Value = '  +1.04%  '

Write-Output "Applied crypto price updates"
